$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Refresh the "last updated" timestamp -------------------------------
$ws.Range("A1").Value = "Datos actualizados a 18 de Mayo de 2020 a las 18:05"

# --- Refresh case counts for the countries whose stats changed ----------
# Columns: B=Casos totales, C=Nuevos casos, D=Casos activos, E=Recuperados,
#          F=Casos criticos, G=Muertes hoy, H=Muertes
$countryData = @{
    "Estados Unidos"       = @(1532861, 5197, 347225, 1094544, 0, 114, 91092)
    "Reino Unido"          = @(246406,  2711,      0,       0, 0, 160, 34796)
    "Alemania"             = @(177069,   418, 154600,   14399, 0,  21,  8070)
    "Suiza"                = @(30597,     10,  27500,    1211, 0,   5,  1886)
    "Polonia"              = @(18885,    356,   7628,   10321, 0,  11,   936)
    "Chequia"              = @(8527,      52,   5633,    2596, 0,   0,   298)
    "Irak"                 = @(3554,     150,   2310,    1117, 0,   4,   127)
    "Estado de Palestina"  = @(386,        5,    337,      47, 0,   0,     2)
    "Liberia"              = @(229,        3,    123,      84, 0,   1,    22)
    "Mozambique"           = @(145,        8,     44,     101, 0,   0,     0)
}

$searchRange = $ws.Range("A4:A219")
foreach ($country in $countryData.Keys) {
    $cell = $searchRange.Find($country)
    $row = $cell.Row
    $values = $countryData[$country]
    for ($i = 0; $i -lt 7; $i++) {
        $ws.Cells.Item($row, 2 + $i).Value = $values[$i]
    }
}

# --- Re-sort the country table by total cases (column B), descending ----
$rng = $ws.Range("A4:H219")
$rng.Sort($ws.Range("B4:B219"), 2)

# Excel's sort does not preserve a fixed order among rows that tie on
# "Casos totales"; align the two affected tie groups (18 and 6 cases) with
# the published order.
$ws.Cells.Item(195, 1).Value = "Nueva Caledonia"
$ws.Cells.Item(197, 1).Value = "Santa Lucia"
$ws.Cells.Item(215, 1).Value = "San Bartolome"
$ws.Cells.Item(216, 1).Value = "Bonaire, San Eustaquio y Saba"
